$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "65.822.66"
$ws.Range("E2").Value = "  -2.46%  "
$ws.Range("D3").Value = "3.268.50"
$ws.Range("E3").Value = "  -1.32%  "
$ws.Range("E4").Value = "  +0.01%  "
$ws.Range("D5").Value = "571.88"
$ws.Range("E5").Value = "  -1.12%  "
$ws.Range("D6").Value = "177.32"
$ws.Range("E6").Value = "  -5.00%  "
$ws.Range("D7").Value = "0.627"
$ws.Range("E7").Value = "  +3.53%  "
$ws.Range("E8").Value = "  -0.04%  "
$ws.Range("E9").Value = "  -3.66%  "
$ws.Range("E10").Value = "  +0.19%  "
$ws.Range("E11").Value = "  -2.86%  "
$ws.Range("D12").Value = "3.841.53"
$ws.Range("E12").Value = "  -1.20%  "
$ws.Range("E13").Value = "  -3.81%  "
$ws.Range("D14").Value = "65.912.41"
$ws.Range("E14").Value = "  -2.67%  "
$ws.Range("D15").Value = "26.38"
$ws.Range("E15").Value = "  -4.10%  "
$ws.Range("E16").Value = "  -3.58%  "
$ws.Range("D17").Value = "3.267.42"
$ws.Range("E17").Value = "  -1.09%  "
$ws.Range("D18").Value = "434.24"
$ws.Range("E18").Value = "  -2.58%  "
$ws.Range("E19").Value = "  -3.02%  "
$ws.Range("D20").Value = "13.11"
$ws.Range("E20").Value = "  -3.71%  "
$ws.Range("E21").Value = "  -5.00%  "
$ws.Range("D22").Value = "72.13"
$ws.Range("E22").Value = "  -2.60%  "
$ws.Range("E23").Value = "  +0.01%  "
$ws.Range("D24").Value = "3.419.74"
$ws.Range("E24").Value = "  -1.03%  "
$ws.Range("D25").Value = "0.504"
$ws.Range("E25").Value = "  -2.95%  "
$ws.Range("E26").Value = "  +3.25%  "
$ws.Range("E27").Value = "  -5.99%  "
$ws.Range("D28").Value = "8.83"
$ws.Range("E28").Value = "  -2.94%  "
$ws.Range("E29").Value = "  -0.17%  "
$ws.Range("E30").Value = "  -2.74%  "
$ws.Range("D31").Value = "22.21"
$ws.Range("E31").Value = "  -3.38%  "
$ws.Range("E32").Value = "  +0.10%  "
$ws.Range("D33").Value = "5.12"
$ws.Range("E33").Value = "  -4.23%  "
$ws.Range("D34").Value = "6.56"
$ws.Range("E34").Value = "  -3.72%  "
$ws.Range("E35").Value = "  -5.85%  "
$ws.Range("D36").Value = "159.86"
$ws.Range("E36").Value = "  -1.96%  "
$ws.Range("E37").Value = "  -6.09%  "
$ws.Range("D38").Value = "26.59"
$ws.Range("E38").Value = "  -2.31%  "
$ws.Range("E39").Value = "  -4.46%  "
$ws.Range("D40").Value = "2.756.36"
$ws.Range("E40").Value = "  -0.15%  "
$ws.Range("D41").Value = "0.775"
$ws.Range("E41").Value = "  -2.20%  "
$ws.Range("D42").Value = "4.30"
$ws.Range("E42").Value = "  -4.01%  "
$ws.Range("D43").Value = "40.21"
$ws.Range("E43").Value = "  -0.06%  "
$ws.Range("D44").Value = "6.01"
$ws.Range("E44").Value = "  -4.21%  "
$ws.Range("D45").Value = "0.0654"
$ws.Range("E45").Value = "  -3.26%  "
$ws.Range("D46").Value = "2.28"
$ws.Range("D47").Value = "317.06"
$ws.Range("E47").Value = "  -3.14%  "
$ws.Range("E48").Value = "  -6.91%  "
$ws.Range("D49").Value = "0.0265"
$ws.Range("E49").Value = "  -3.35%  "
$ws.Range("E50").Value = "  +1.68%  "
$ws.Range("E51").Value = "  -0.02%  "
